$wb = $excel.ActiveWorkbook

# --- Populate the "ScheduleData" sheet with the new schedule table ---
$ws = $wb.Worksheets.Item("ScheduleData")

# Header row, left to right
$ws.Range("A1").Value = "ROT"
$ws.Range("B1").Value = "ID"
$ws.Range("C1").Value = "First"
$ws.Range("D1").Value = "Last"
$ws.Range("E1").Value = "ROT1"
$ws.Range("F1").Value = "ROT2"
$ws.Range("G1").Value = "ROT3"
$ws.Range("H1").Value = "ROT4"

# Row 2 rotation values (station labels first), then names
$ws.Range("F2").Value = "H"
$ws.Range("G2").Value = "S"
$ws.Range("E2").Value = "GE"
$ws.Range("H2").Value = "GL"
$ws.Range("C2").Value = "Tom"
$ws.Range("D2").Value = "Tim"

# Row 3 names, then rotation values (rotated one station over)
$ws.Range("C3").Value = "Jim"
$ws.Range("D3").Value = "John"
$ws.Range("E3").Value = "H"
$ws.Range("F3").Value = "S"
$ws.Range("G3").Value = "GL"
$ws.Range("H3").Value = "GE"

# Select the range that was last selected on this sheet, then make it active
$ws.Range("C3:H3").Select()
$ws.Activate()

# --- Update selection on the "RawData" sheet ---
$ws1 = $wb.Worksheets.Item("RawData")
$ws1.Range("F12").Select()

# Re-activate ScheduleData so it ends up as the selected/active tab
$ws.Activate()
